# Update excess mortality analyses

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berekening oversterfte")

# --- Updated weekly "waargenomen" (observed) figures for existing weeks ---
$ws.Range("G17").Value = 2694
$ws.Range("G23").Value = 2659
$ws.Range("G24").Value = 2638
$ws.Range("G25").Value = 3206
$ws.Range("G26").Value = 2845
$ws.Range("G27").Value = 2731
$ws.Range("G29").Value = 2733
$ws.Range("G30").Value = 2710
$ws.Range("G31").Value = 2881
$ws.Range("G32").Value = 2988
$ws.Range("G33").Value = 2996

# --- New week 42 row (row 34) ---
$ws.Range("F34").Value = 42
$ws.Range("G34").Value = 3224
$ws.Range("H34").Value = 2839
$ws.Range("I34").Formula = "=G34-H34"

# --- Extend the year total to include the new week ---
$ws.Range("I37").Formula = "=SUM(I3:I34)"

# --- Update view state: selection moved to I38 (scrolled so row 2 is at top) ---
$ws.Range("I38").Select()

$wb.Save()
